$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove " with interest" from the opening sentence, splitting the
#    original single run into two runs around the comma (the _GoBack
#    bookmark will be relocated to that split point below).
# ------------------------------------------------------------------
$d.Content.Find.Execute(" with interest", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# ------------------------------------------------------------------
# 2) Relocate the "_GoBack" bookmark from its old position (between
#    "veflow" and " profile") to right after "...Big Data Engineer"
#    (i.e. right before the comma that now starts the second run).
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$r = $d.Content
$r.Find.Execute("As a highly skilled Big Data Engineer", $true, $false, $false, $false, `
                $false, $true, 1, $false, "", 0)
$splitPoint = $r.End
$newRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $newRange)

# ------------------------------------------------------------------
# 3) Merge the "veflow" and " profile" runs into a single run
#    ("veflow profile") now that the bookmark no longer separates
#    them, while leaving the surrounding "O" run and the trailing
#    " for additional details..." run untouched. We do this by
#    temporarily nudging the formatting of those two neighbours so
#    the engine's run-coalescing won't fold them into the merge, then
#    restoring them afterwards (a no-op, byte-identical restore).
# ------------------------------------------------------------------
$start = $d.Content
$start.Find.Execute("resume and S", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterS = $d.Range($start.End, $d.Content.End)
$afterS.Find.Execute("tack", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterTack = $d.Range($afterS.End, $d.Content.End)
$afterTack.Find.Execute("O", $true, $false, $false, $false, $true, $true, 1, $false, "", 0)
$oRun = $d.Range($afterTack.Start, $afterTack.End)

$afterO = $d.Range($oRun.End, $d.Content.End)
$afterO.Find.Execute("veflow", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$veflowRun = $d.Range($afterO.Start, $afterO.End)

$afterVeflow = $d.Range($veflowRun.End, $d.Content.End)
$afterVeflow.Find.Execute(" profile", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$profileRun = $d.Range($afterVeflow.Start, $afterVeflow.End)

$afterProfile = $d.Range($profileRun.End, $d.Content.End)
$afterProfile.Find.Execute(" for additional details regarding my expertise and career achievements. I am glad to discuss with you how my experience and background meet your needs.", `
                            $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$forRun = $d.Range($afterProfile.Start, $afterProfile.End)

$oRun.Bold = 1
$forRun.Bold = 1

$mergeSpan = $d.Range($veflowRun.Start, $profileRun.End)
$mergeSpan.Find.Execute("veflow profile", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "veflow profile", 1)

$oRun.Bold = 0
$forRun.Bold = 0
